# Added new employee in PIM module
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# New header columns (match formatting of the existing header row)
$ws.Range("K1").Value = "FirstName"
$ws.Range("L1").Value = "LastName"
$ws.Range("M1").Value = "EmployeeId"
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1:M1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Existing filter rows (2-6) get an empty quote-prefixed text value in the
# new columns, matching the rest of the table's "blank placeholder" cells.
$ws.Range("K2").Value = "'"
$ws.Range("L2").Value = "'"
$ws.Range("M2").Value = "'"

$ws.Range("K3").Value = "'"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'"

$ws.Range("K4").Value = "'"
$ws.Range("L4").Value = "'"
$ws.Range("M4").Value = "'"

$ws.Range("K5").Value = "'"
$ws.Range("L5").Value = "'"
$ws.Range("M5").Value = "'"

$ws.Range("K6").Value = "'"
$ws.Range("L6").Value = "'"
$ws.Range("M6").Value = "'"

# New employee test data row
$ws.Range("K7").Value = "neelesh"
$ws.Range("L7").Value = "khatri"
$ws.Range("M7").Value = "neel"

# Approximate the auto-fit column widths for the new columns
$ws.Columns.Item(11).ColumnWidth = 8.6
$ws.Columns.Item(12).ColumnWidth = 8.5
$ws.Columns.Item(13).ColumnWidth = 10.166666666666666

# Scroll the view over and select the newly added cell
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
$null = $ws.Range("M7").Select()
